$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair data (natmi analysis, following Dr Hou advice) for rows 2-11, columns A:T
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
#          E..T = numeric NATMI metrics
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @("ECs", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 51.74808266666667, 155.244248, 0.1468432527076591, 0.1468432527076591, 2, 0.6666666666666666, 0.2083573333333334, 0.6250720000000001, 0.479707634607725, 0.4797076346077249, 10.78209250953956, 97.03883258585601, 0.07044182941449556, 0.07044182941449555),
    @("ECs", "Rtn4", "Tnfrsf19", "sCs", 3, 1, 51.74808266666667, 155.244248, 0.1468432527076591, 0.1468432527076591, 3, 1, 0.225985, 0.6779550000000001, 0.5202923653922751, 0.5202923653922751, 11.69429046142667, 105.24861415284, 0.07640142329316356, 0.07640142329316356),
    @("FAPs", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 64.52529266666666, 193.575878, 0.1831005782015253, 0.1831005782015253, 2, 0.6666666666666666, 0.2083573333333334, 0.6250720000000001, 0.479707634607725, 0.4797076346077249, 13.44431791257956, 120.998861213216, 0.08783474526436046, 0.08783474526436047),
    @("FAPs", "Rtn4", "Tnfrsf19", "sCs", 3, 1, 64.52529266666666, 193.575878, 0.1831005782015253, 0.1831005782015253, 3, 1, 0.225985, 0.6779550000000001, 0.5202923653922751, 0.5202923653922751, 14.58174826327667, 131.23573436949, 0.09526583293716483, 0.09526583293716484),
    @("M1", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 102.0393243333333, 306.117973, 0.2895524919389952, 0.2895524919389952, 2, 0.6666666666666666, 0.2083573333333334, 0.6250720000000001, 0.479707634607725, 0.4797076346077249, 21.26064151322845, 191.345773619056, 0.1389005410028277, 0.1389005410028277),
    @("M1", "Rtn4", "Tnfrsf19", "sCs", 3, 1, 102.0393243333333, 306.117973, 0.2895524919389952, 0.2895524919389952, 3, 1, 0.225985, 0.6779550000000001, 0.5202923653922751, 0.5202923653922751, 23.05935670946834, 207.534210385215, 0.1506519509361675, 0.1506519509361675),
    @("M2", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 84.450424, 253.351272, 0.239641244924597, 0.239641244924597, 2, 0.6666666666666666, 0.2083573333333334, 0.6250720000000001, 0.479707634607725, 0.4797076346077249, 17.59586514350934, 158.362786291584, 0.1149577347572289, 0.1149577347572289),
    @("M2", "Rtn4", "Tnfrsf19", "sCs", 3, 1, 84.450424, 253.351272, 0.239641244924597, 0.239641244924597, 3, 1, 0.225985, 0.6779550000000001, 0.5202923653922751, 0.5202923653922751, 19.08452906764, 171.76076160876, 0.1246835101673681, 0.1246835101673681),
    @("sCs", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 49.64042033333334, 148.921261, 0.1408624322272234, 0.1408624322272234, 2, 0.6666666666666666, 0.2083573333333334, 0.6250720000000001, 0.479707634607725, 0.4797076346077249, 10.34294560619911, 93.08651045579202, 0.06757278416881229, 0.06757278416881229),
    @("sCs", "Rtn4", "Tnfrsf19", "sCs", 3, 1, 49.64042033333334, 148.921261, 0.1408624322272234, 0.1408624322272234, 3, 1, 0.225985, 0.6779550000000001, 0.5202923653922751, 0.5202923653922751, 11.21799038902834, 100.961913501255, 0.07328964805841109, 0.07328964805841109),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
